$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 content changes ---
# H4: "Underway" -> TRUE
$ws1.Range("H4").Value = $true
# H12: "Underway" -> TRUE
$ws1.Range("H12").Value = $true
# G13: FALSE -> TRUE
$ws1.Range("G13").Value = $true
# L13: "Waiting for me" -> "Yes"
$ws1.Range("L13").Value = "Yes"

# --- Sheet2 content changes ---
# H2: empty -> "RabbiBordon"
$ws2.Range("H2").Value = "RabbiBordon"
# D6: apply "Good" cell style (green highlight)
$ws2.Range("D6").Style = "Good"

# --- Selection / active sheet changes ---
$ws2.Activate()
$ws2.Range("C7").Select()

$ws1.Activate()
$ws1.Range("B5").Select()
